$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell (AC1) onto the
# new header cells AD1:AF1, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill data rows 2-50 with team record: Wins=73, Losses=89, Ties=0
$lastRow = 50
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 73   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 89   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
